# The "퇴사율" (turnover-rate) sheet had its "전체" (total) column removed.
# Deleting the whole column shifts the remaining columns left (C:H -> B:G),
# which in turn drops the now-unused "전체" shared string and makes this
# sheet the active tab/selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("퇴사율")
$ws.Activate()
$ws.Columns("B").Delete()
$ws.Range("C12").Select()
